# Presentation state 11.02 - naive component forecaster fix.
# Every forecast-vintage row gains a new leading "next quarter" (Q0)
# error value; the values that used to occupy columns B:J shift one
# column to the right (C:K), and the oldest trailing quarter (which
# falls off the 10-quarter window) is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the oldest vintage has no room left for a new leading quarter,
# so it simply drops its stale trailing value (K2).
$ws.Cells.Item(2, 11).ClearContents()

# New Q0 error value introduced at column B for each vintage row.
$newB = @{
    3 = 2.174397541324862
    4 = 8.987815262142332
    5 = -8.470544361886736
    6 = -0.34111357821662
    7 = 1.338580889567239
    8 = -1.631436389465022
    9 = -1.810267219091221
    10 = 0.7496711949059137
    11 = -0.1976049264540507
    12 = 0.1255103924969555
    13 = -0.003050974015260888
    14 = 1.419652253737389
    15 = 0.1883110177716323
    16 = 0.2202779152847414
    17 = 0.4485660054549828
    18 = 0.420735823599318
    19 = -0.1252553916527783
    20 = 0.09764018641116785
    21 = -0.1898380159455487
    22 = 0.1743923273248104
    23 = -0.254916590923889
    24 = -0.01954473626955232
}

# Last populated column (B=2 .. K=11) in each row before the edit.
$lastColBefore = @{
    3 = 11
    4 = 11
    5 = 11
    6 = 11
    7 = 11
    8 = 11
    9 = 11
    10 = 11
    11 = 11
    12 = 11
    13 = 11
    14 = 11
    15 = 10
    16 = 9
    17 = 8
    18 = 7
    19 = 6
    20 = 5
    21 = 4
    22 = 3
    23 = 2
    24 = 1
}

# Rows 3-24: shift existing values right by one column (right-to-left so
# a value is never overwritten before it has been read), then drop in the
# new leading Q0 value at column B. Column K (11) is the last column in
# the sheet, so a row that was already full out to K (10 values, B:K) can
# only shift its B:J values into C:K -- the old K value (the oldest, 10th
# quarter) simply falls off the window instead of moving to a column L.
foreach ($row in 3..24) {
    $last = $lastColBefore[$row]
    $shiftFrom = [Math]::Min($last, 10)
    for ($col = $shiftFrom; $col -ge 2; $col--) {
        $src = $ws.Cells.Item($row, $col).Value2
        $ws.Cells.Item($row, $col + 1).Value2 = $src
    }
    $ws.Cells.Item($row, 2).Value2 = $newB[$row]
}
